$d = $word.ActiveDocument

$d.Content.Find.Execute("25+62=", $true, $false, $false, $false, $false, $true, 1, $false, "47+39=", 2) | Out-Null
$d.Content.Find.Execute("47+6=", $true, $false, $false, $false, $false, $true, 1, $false, "1+90=", 2) | Out-Null
$d.Content.Find.Execute("70-18=", $true, $false, $false, $false, $false, $true, 1, $false, "97-7=", 2) | Out-Null
$d.Content.Find.Execute("30+0=", $true, $false, $false, $false, $false, $true, 1, $false, "25-4=", 2) | Out-Null
$d.Content.Find.Execute("80-44=", $true, $false, $false, $false, $false, $true, 1, $false, "90-56=", 2) | Out-Null
$d.Content.Find.Execute("89-16=", $true, $false, $false, $false, $false, $true, 1, $false, "83-48=", 2) | Out-Null
$d.Content.Find.Execute("41+48=", $true, $false, $false, $false, $false, $true, 1, $false, "5+35=", 2) | Out-Null
$d.Content.Find.Execute("27+14=", $true, $false, $false, $false, $false, $true, 1, $false, "13+38=", 2) | Out-Null
$d.Content.Find.Execute("58-11=", $true, $false, $false, $false, $false, $true, 1, $false, "47-10=", 2) | Out-Null
$d.Content.Find.Execute("46-21=", $true, $false, $false, $false, $false, $true, 1, $false, "96-83=", 2) | Out-Null
$d.Content.Find.Execute("49-5=", $true, $false, $false, $false, $false, $true, 1, $false, "1+66=", 2) | Out-Null
$d.Content.Find.Execute("95-28=", $true, $false, $false, $false, $false, $true, 1, $false, "94+1=", 2) | Out-Null
$d.Content.Find.Execute("62-57=", $true, $false, $false, $false, $false, $true, 1, $false, "51+24=", 2) | Out-Null
$d.Content.Find.Execute("40-24=", $true, $false, $false, $false, $false, $true, 1, $false, "5+68=", 2) | Out-Null
$d.Content.Find.Execute("0+41=", $true, $false, $false, $false, $false, $true, 1, $false, "18+45=", 2) | Out-Null
$d.Content.Find.Execute("37+19=", $true, $false, $false, $false, $false, $true, 1, $false, "78-3=", 2) | Out-Null
$d.Content.Find.Execute("45+49=", $true, $false, $false, $false, $false, $true, 1, $false, "60+11=", 2) | Out-Null
$d.Content.Find.Execute("33-10=", $true, $false, $false, $false, $false, $true, 1, $false, "3+12=", 2) | Out-Null
$d.Content.Find.Execute("12+85=", $true, $false, $false, $false, $false, $true, 1, $false, "65-39=", 2) | Out-Null
$d.Content.Find.Execute("87-17=", $true, $false, $false, $false, $false, $true, 1, $false, "27+30=", 2) | Out-Null
$d.Content.Find.Execute("94-94=", $true, $false, $false, $false, $false, $true, 1, $false, "68-45=", 2) | Out-Null
$d.Content.Find.Execute("21-9=", $true, $false, $false, $false, $false, $true, 1, $false, "9-5=", 2) | Out-Null
$d.Content.Find.Execute("4+86=", $true, $false, $false, $false, $false, $true, 1, $false, "47-19=", 2) | Out-Null
$d.Content.Find.Execute("40+53=", $true, $false, $false, $false, $false, $true, 1, $false, "38-3=", 2) | Out-Null
$d.Content.Find.Execute("1+44=", $true, $false, $false, $false, $false, $true, 1, $false, "52-49=", 2) | Out-Null
$d.Content.Find.Execute("91-85=", $true, $false, $false, $false, $false, $true, 1, $false, "24+12=", 2) | Out-Null
$d.Content.Find.Execute("13-6=", $true, $false, $false, $false, $false, $true, 1, $false, "82-21=", 2) | Out-Null
$d.Content.Find.Execute("69-7=", $true, $false, $false, $false, $false, $true, 1, $false, "55-53=", 2) | Out-Null
$d.Content.Find.Execute("90-12=", $true, $false, $false, $false, $false, $true, 1, $false, "9-1=", 2) | Out-Null
$d.Content.Find.Execute("30+55=", $true, $false, $false, $false, $false, $true, 1, $false, "89-48=", 2) | Out-Null
$d.Content.Find.Execute("29+27=", $true, $false, $false, $false, $false, $true, 1, $false, "20+53=", 2) | Out-Null
$d.Content.Find.Execute("53+13=", $true, $false, $false, $false, $false, $true, 1, $false, "10+30=", 2) | Out-Null
$d.Content.Find.Execute("88-54=", $true, $false, $false, $false, $false, $true, 1, $false, "72-45=", 2) | Out-Null
$d.Content.Find.Execute("58+25=", $true, $false, $false, $false, $false, $true, 1, $false, "11+68=", 2) | Out-Null
$d.Content.Find.Execute("82-77=", $true, $false, $false, $false, $false, $true, 1, $false, "75-12=", 2) | Out-Null
$d.Content.Find.Execute("18+78=", $true, $false, $false, $false, $false, $true, 1, $false, "97-21=", 2) | Out-Null
$d.Content.Find.Execute("22-0=", $true, $false, $false, $false, $false, $true, 1, $false, "50+4=", 2) | Out-Null
$d.Content.Find.Execute("21+3=", $true, $false, $false, $false, $false, $true, 1, $false, "26+43=", 2) | Out-Null
$d.Content.Find.Execute("1+50=", $true, $false, $false, $false, $false, $true, 1, $false, "46+2=", 2) | Out-Null
$d.Content.Find.Execute("56-46=", $true, $false, $false, $false, $false, $true, 1, $false, "86+13=", 2) | Out-Null
$d.Content.Find.Execute("95-61=", $true, $false, $false, $false, $false, $true, 1, $false, "15-8=", 2) | Out-Null
$d.Content.Find.Execute("29-29=", $true, $false, $false, $false, $false, $true, 1, $false, "90+5=", 2) | Out-Null
$d.Content.Find.Execute("90-41=", $true, $false, $false, $false, $false, $true, 1, $false, "26+58=", 2) | Out-Null
$d.Content.Find.Execute("69+20=", $true, $false, $false, $false, $false, $true, 1, $false, "56-22=", 2) | Out-Null
$d.Content.Find.Execute("10+14=", $true, $false, $false, $false, $false, $true, 1, $false, "26+69=", 2) | Out-Null
$d.Content.Find.Execute("89-20=", $true, $false, $false, $false, $false, $true, 1, $false, "48-5=", 2) | Out-Null
$d.Content.Find.Execute("94-15=", $true, $false, $false, $false, $false, $true, 1, $false, "77-71=", 2) | Out-Null
$d.Content.Find.Execute("38-27=", $true, $false, $false, $false, $false, $true, 1, $false, "79-26=", 2) | Out-Null
$d.Content.Find.Execute("69-68=", $true, $false, $false, $false, $false, $true, 1, $false, "52+39=", 2) | Out-Null
$d.Content.Find.Execute("59-22=", $true, $false, $false, $false, $false, $true, 1, $false, "20+20=", 2) | Out-Null
$d.Content.Find.Execute("2+19=", $true, $false, $false, $false, $false, $true, 1, $false, "24+30=", 2) | Out-Null
$d.Content.Find.Execute("90-79=", $true, $false, $false, $false, $false, $true, 1, $false, "28-2=", 2) | Out-Null
$d.Content.Find.Execute("35+3=", $true, $false, $false, $false, $false, $true, 1, $false, "7+76=", 2) | Out-Null
$d.Content.Find.Execute("79+10=", $true, $false, $false, $false, $false, $true, 1, $false, "27+40=", 2) | Out-Null
$d.Content.Find.Execute("27+32=", $true, $false, $false, $false, $false, $true, 1, $false, "8+40=", 2) | Out-Null
$d.Content.Find.Execute("98-70=", $true, $false, $false, $false, $false, $true, 1, $false, "25-15=", 2) | Out-Null
$d.Content.Find.Execute("56+15=", $true, $false, $false, $false, $false, $true, 1, $false, "4+21=", 2) | Out-Null
$d.Content.Find.Execute("91-50=", $true, $false, $false, $false, $false, $true, 1, $false, "75+4=", 2) | Out-Null
$d.Content.Find.Execute("48+25=", $true, $false, $false, $false, $false, $true, 1, $false, "86-21=", 2) | Out-Null
$d.Content.Find.Execute("89-31=", $true, $false, $false, $false, $false, $true, 1, $false, "70-37=", 2) | Out-Null
$d.Content.Find.Execute("18+46=", $true, $false, $false, $false, $false, $true, 1, $false, "51-1=", 2) | Out-Null
$d.Content.Find.Execute("66+29=", $true, $false, $false, $false, $false, $true, 1, $false, "53+35=", 2) | Out-Null
$d.Content.Find.Execute("12+78=", $true, $false, $false, $false, $false, $true, 1, $false, "75-52=", 2) | Out-Null
$d.Content.Find.Execute("24+51=", $true, $false, $false, $false, $false, $true, 1, $false, "10+6=", 2) | Out-Null
$d.Content.Find.Execute("12+29=", $true, $false, $false, $false, $false, $true, 1, $false, "23+25=", 2) | Out-Null
$d.Content.Find.Execute("24-9=", $true, $false, $false, $false, $false, $true, 1, $false, "20+25=", 2) | Out-Null
$d.Content.Find.Execute("45+37=", $true, $false, $false, $false, $false, $true, 1, $false, "89-32=", 2) | Out-Null
$d.Content.Find.Execute("12+37=", $true, $false, $false, $false, $false, $true, 1, $false, "8+69=", 2) | Out-Null
$d.Content.Find.Execute("17-1=", $true, $false, $false, $false, $false, $true, 1, $false, "32+35=", 2) | Out-Null
$d.Content.Find.Execute("16+21=", $true, $false, $false, $false, $false, $true, 1, $false, "63+23=", 2) | Out-Null
$d.Content.Find.Execute("15+10=", $true, $false, $false, $false, $false, $true, 1, $false, "13+12=", 2) | Out-Null
$d.Content.Find.Execute("59-42=", $true, $false, $false, $false, $false, $true, 1, $false, "8+77=", 2) | Out-Null
$d.Content.Find.Execute("17-0=", $true, $false, $false, $false, $false, $true, 1, $false, "55+32=", 2) | Out-Null
$d.Content.Find.Execute("93-20=", $true, $false, $false, $false, $false, $true, 1, $false, "89-46=", 2) | Out-Null
$d.Content.Find.Execute("68-7=", $true, $false, $false, $false, $false, $true, 1, $false, "84-63=", 2) | Out-Null
$d.Content.Find.Execute("48-9=", $true, $false, $false, $false, $false, $true, 1, $false, "0+35=", 2) | Out-Null
$d.Content.Find.Execute("62-50=", $true, $false, $false, $false, $false, $true, 1, $false, "39+26=", 2) | Out-Null
$d.Content.Find.Execute("31+36=", $true, $false, $false, $false, $false, $true, 1, $false, "84-26=", 2) | Out-Null
$d.Content.Find.Execute("27+59=", $true, $false, $false, $false, $false, $true, 1, $false, "8+75=", 2) | Out-Null
$d.Content.Find.Execute("11+21=", $true, $false, $false, $false, $false, $true, 1, $false, "3+80=", 2) | Out-Null
$d.Content.Find.Execute("49-32=", $true, $false, $false, $false, $false, $true, 1, $false, "2+65=", 2) | Out-Null
$d.Content.Find.Execute("48+43=", $true, $false, $false, $false, $false, $true, 1, $false, "3+83=", 2) | Out-Null
$d.Content.Find.Execute("88-65=", $true, $false, $false, $false, $false, $true, 1, $false, "19-19=", 2) | Out-Null
$d.Content.Find.Execute("7+36=", $true, $false, $false, $false, $false, $true, 1, $false, "60-26=", 2) | Out-Null
$d.Content.Find.Execute("88-80=", $true, $false, $false, $false, $false, $true, 1, $false, "85-44=", 2) | Out-Null
$d.Content.Find.Execute("84-78=", $true, $false, $false, $false, $false, $true, 1, $false, "28+21=", 2) | Out-Null
$d.Content.Find.Execute("0+14=", $true, $false, $false, $false, $false, $true, 1, $false, "17+38=", 2) | Out-Null
$d.Content.Find.Execute("15+8=", $true, $false, $false, $false, $false, $true, 1, $false, "66-12=", 2) | Out-Null
$d.Content.Find.Execute("39+46=", $true, $false, $false, $false, $false, $true, 1, $false, "99-53=", 2) | Out-Null
$d.Content.Find.Execute("47+2=", $true, $false, $false, $false, $false, $true, 1, $false, "47+31=", 2) | Out-Null
$d.Content.Find.Execute("28-0=", $true, $false, $false, $false, $false, $true, 1, $false, "99-99=", 2) | Out-Null
$d.Content.Find.Execute("42-41=", $true, $false, $false, $false, $false, $true, 1, $false, "80-23=", 2) | Out-Null
$d.Content.Find.Execute("76+14=", $true, $false, $false, $false, $false, $true, 1, $false, "8+82=", 2) | Out-Null
$d.Content.Find.Execute("94-0=", $true, $false, $false, $false, $false, $true, 1, $false, "22+63=", 2) | Out-Null
$d.Content.Find.Execute("63+16=", $true, $false, $false, $false, $false, $true, 1, $false, "23-9=", 2) | Out-Null
$d.Content.Find.Execute("81+3=", $true, $false, $false, $false, $false, $true, 1, $false, "65-55=", 2) | Out-Null
$d.Content.Find.Execute("13+82=", $true, $false, $false, $false, $false, $true, 1, $false, "14-1=", 2) | Out-Null
$d.Content.Find.Execute("79-65=", $true, $false, $false, $false, $false, $true, 1, $false, "30+27=", 2) | Out-Null
$d.Content.Find.Execute("10-9=", $true, $false, $false, $false, $false, $true, 1, $false, "35-11=", 2) | Out-Null
$d.Content.Find.Execute("94-6=", $true, $false, $false, $false, $false, $true, 1, $false, "32+5=", 2) | Out-Null
Write-Host "Replacements complete"
